$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Typography")

# Row 6 ("Small" typography) prototype values, per latest UI draft.
$ws.Range("F6").Value = "_"
$ws.Range("G6").Value = " &é`"'(-è_çà)=~#{[|``\^@]}°+>,;:!?./§ù*^$%µ£¨"
$ws.Range("I6").Value = "a-z,A-Z,0-9"
